# ore_tirocinio.docx — registra la seconda giornata di tirocinio
# (creazione script per generazione gif NDVI Campania) e aggiorna l'anno
# della prima data da un refuso 2022 -> 2023.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Ridimensiona le colonne della tabella (tblGrid) -------------------
$t.Columns.Item(1).Width = 101.2   # 2024 dxa
$t.Columns.Item(2).Width = 86.45   # 1729 dxa
$t.Columns.Item(3).Width = 76.5    # 1530 dxa
$t.Columns.Item(4).Width = 217.25  # 4345 dxa

# --- 2. Corregge l'anno della prima data -----------------------------------
$d.Content.Find.Execute("02/10/2022", $true, $false, $false, $false, $false, `
    $true, 1, $false, "02/10/2023", 2) | Out-Null

# --- 3. Unisce i due run "11" + ":00 " della prima riga in un unico run ----
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$cellOrario1 = $t.Cell(2, 2)
$xmlOrario1 = '<w:p ' + $ns + '>' + `
    '<w:pPr><w:jc w:val="center"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">11:00 </w:t></w:r>' + `
    '<w:r><w:t>' + [char]0x2013 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>16:00</w:t></w:r>' + `
    '</w:p>'
$cellOrario1.Range.InsertXML($xmlOrario1)

# --- 4. Rimuove il bookmark _GoBack dalla fine della prima riga ------------
$cellAttivita1 = $t.Cell(2, 4)
$xmlAttivita1 = '<w:p ' + $ns + '>' + `
    '<w:r><w:t xml:space="preserve">Settaggio ambiente di lavoro (Google Earth Engine), importazione degli </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>shapefile</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> e creazione script ' + [char]0x201C + 'VegetazioneCampania.js' + [char]0x201D + '.</w:t></w:r>' + `
    '</w:p>'
$cellAttivita1.Range.InsertXML($xmlAttivita1)

# --- 5. Compila la seconda riga della tabella ------------------------------
$cellGiorno2 = $t.Cell(3, 1)
$xmlGiorno2 = '<w:p ' + $ns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
    '<w:r><w:t>03/10/2023</w:t></w:r></w:p>'
$cellGiorno2.Range.InsertXML($xmlGiorno2)

$cellOrario2 = $t.Cell(3, 2)
$xmlOrario2 = '<w:p ' + $ns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
    '<w:r><w:t>10:00 ' + [char]0x2013 + ' 16:00</w:t></w:r></w:p>'
$cellOrario2.Range.InsertXML($xmlOrario2)

$cellOre2 = $t.Cell(3, 3)
$xmlOre2 = '<w:p ' + $ns + '><w:pPr><w:jc w:val="center"/></w:pPr>' + `
    '<w:r><w:t>6</w:t></w:r></w:p>'
$cellOre2.Range.InsertXML($xmlOre2)

$cellAttivita2 = $t.Cell(3, 4)
$quote = [char]0x201C
$quoteEnd = [char]0x201D
$xmlAttivita2 = '<w:p ' + $ns + '>' + `
    '<w:r><w:t>Lettura documentazione ' + $quote + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">MODIS NDVI Times Series </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Animation</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>' + $quoteEnd + ' da Google Earth Engine e creazione script ' + $quote + '</w:t></w:r>' + `
    '<w:r><w:t>generatoreGifVegetazioneCampania2022</w:t></w:r>' + `
    '<w:r><w:t>.js' + $quoteEnd + '.</w:t></w:r>' + `
    '</w:p>'
$cellAttivita2.Range.InsertXML($xmlAttivita2)

# --- 6. Sposta il bookmark _GoBack sul paragrafo finale del documento ------
# (Si usa Range invece di Paragraphs, la cui collezione non si riallinea in
# modo affidabile dopo le InsertXML precedenti.)
$tableEnd = $t.Range.End
$finalPara = $d.Range($tableEnd, $d.Content.End)
$finalPara.InsertXML('<w:p ' + $ns + '>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>')
